{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Change 1: the centered date paragraph \"9/11/2014\" becomes two runs,\n//   \"9/26\" and \"/2014\" (i.e. the date text becomes \"9/26/2014\" but is\n//   represented as two separate runs because only the \"9/11\" portion was\n//   retyped).\n// Change 2: a new paragraph is appended right after the paragraph that\n//   ends in \"...IDmsMetaData interface.\" (and right before the paragraph\n//   that begins \"After compiling your DMS library...\").\n\nconst FLAT_OPC_HEADER =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>';\nconst FLAT_OPC_FOOTER =\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrapFlatOpc(innerBodyXml) {\n  return FLAT_OPC_HEADER + innerBodyXml + FLAT_OPC_FOOTER;\n}\n\nconst body = context.document.body;\n\n// --- Change 1: split \"9/11/2014\" into \"9/26\" + \"/2014\" -----------------\nconst dateResults = body.search(\"9/11/2014\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  const dateRange = dateResults.items[0];\n  const dateOoxml = wrapFlatOpc(\n    \"<w:p><w:r><w:t>9/26</w:t></w:r><w:r><w:t>/2014</w:t></w:r></w:p>\"\n  );\n  dateRange.insertOoxml(dateOoxml, Word.InsertLocation.replace);\n}\n\n// --- Change 2: insert new paragraph after the \"...interface.\" paragraph -\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"The metadata available to the DMS functionality\") >= 0 &&\n      p.text.indexOf(\"interface.\") >= 0) {\n    targetParagraph = p;\n    break;\n  }\n}\n\nif (targetParagraph) {\n  const newParagraphXml =\n    \"<w:p>\" +\n    \"<w:r><w:t>Further, you must create a validator that expo</w:t></w:r>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\">rts the </w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n    \"<w:r><w:t>IDmsValidator</w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\"> interface, and additionally, a </w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n    \"<w:r><w:t>DMSValidation</w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\"> control that inherits from the </w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n    \"<w:r><w:t>classDMSBaseControl</w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\"> class</w:t></w:r>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\">, the type of which will be returned by the validators </w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n    \"<w:r><w:t>DMSValidatorControl</w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\"> property</w:t></w:r>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\">. </w:t></w:r>\" +\n    \"</w:p>\";\n\n  const endOfTarget = targetParagraph.getRange(\"End\");\n  endOfTarget.insertOoxml(wrapFlatOpc(newParagraphXml), Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is already open as $d.\n#\n# Change 1: the centered date paragraph \"9/11/2014\" becomes two runs,\n#   \"9/26\" and \"/2014\" (i.e. the date text becomes \"9/26/2014\" but is\n#   represented as two separate runs because only the \"9/11\" portion was\n#   retyped).\n# Change 2: a new paragraph is appended right after the paragraph that\n#   ends in \"...IDmsMetaData interface.\" (and right before the paragraph\n#   that begins \"After compiling your DMS library...\").\n\n$d = $word.ActiveDocument\n$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'\n\n# --- Change 1: split \"9/11/2014\" into \"9/26\" + \"/2014\" ----------------\n$dateParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text.TrimEnd(\"`r\") -eq \"9/11/2014\") {\n        $dateParagraph = $d.Paragraphs($i)\n        break\n    }\n}\n\nif ($dateParagraph -ne $null) {\n    $dateRange = $dateParagraph.Range\n    $dateXml = \"<w:p xmlns:w='$wNs'><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>9/26</w:t></w:r><w:r><w:t>/2014</w:t></w:r></w:p>\"\n    $dateRange.InsertXML($dateXml) | Out-Null\n}\n\n# --- Change 2: insert new paragraph after the \"...interface.\" paragraph\n$targetParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if (($t -like \"*The metadata available to the DMS functionality*\") -and ($t -like \"*interface.*\")) {\n        $targetParagraph = $d.Paragraphs($i)\n        break\n    }\n}\n\nif ($targetParagraph -ne $null) {\n    $targetParagraph.Range.InsertParagraphAfter()\n    $newIndex = $targetParagraph.Index + 1\n    $newParaRange = $d.Paragraphs($newIndex).Range\n\n    $innerXml = '<w:r><w:t>Further, you must create a validator that expo</w:t></w:r>' + `\n        '<w:r><w:t xml:space=\"preserve\">rts the </w:t></w:r>' + `\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>IDmsValidator</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' + `\n        '<w:r><w:t xml:space=\"preserve\"> interface, and additionally, a </w:t></w:r>' + `\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>DMSValidation</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' + `\n        '<w:r><w:t xml:space=\"preserve\"> control that inherits from the </w:t></w:r>' + `\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>classDMSBaseControl</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' + `\n        '<w:r><w:t xml:space=\"preserve\"> class</w:t></w:r>' + `\n        '<w:r><w:t xml:space=\"preserve\">, the type of which will be returned by the validators </w:t></w:r>' + `\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>DMSValidatorControl</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' + `\n        '<w:r><w:t xml:space=\"preserve\"> property</w:t></w:r>' + `\n        '<w:r><w:t xml:space=\"preserve\">. </w:t></w:r>'\n\n    $newParaXml = \"<w:p xmlns:w='$wNs'>$innerXml</w:p>\"\n    $newParaRange.InsertXML($newParaXml) | Out-Null\n}\n"}
